$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-13 (Wnt9a-Fzd4 pairs) with refreshed TPM-derived values
# and drop the "Resolving-Mac" target-cluster rows that no longer exist in the new run.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt9a"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.217872
$ws.Cells.Item(2, 8).Value = 3.653616
$ws.Cells.Item(2, 9).Value = 0.4812552202382371
$ws.Cells.Item(2, 10).Value = 0.4812552202382372
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 26.60444266666667
$ws.Cells.Item(2, 14).Value = 79.813328
$ws.Cells.Item(2, 15).Value = 0.5736225649467147
$ws.Cells.Item(2, 16).Value = 0.5736225649467147
$ws.Cells.Item(2, 17).Value = 32.40080579933867
$ws.Cells.Item(2, 18).Value = 291.607252194048
$ws.Cells.Item(2, 19).Value = 0.2760588538270537
$ws.Cells.Item(2, 20).Value = 0.2760588538270537

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt9a"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.217872
$ws.Cells.Item(3, 8).Value = 3.653616
$ws.Cells.Item(3, 9).Value = 0.4812552202382371
$ws.Cells.Item(3, 10).Value = 0.4812552202382372
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 14.03147833333333
$ws.Cells.Item(3, 14).Value = 42.094435
$ws.Cells.Item(3, 15).Value = 0.3025349071358453
$ws.Cells.Item(3, 16).Value = 0.3025349071358453
$ws.Cells.Item(3, 17).Value = 17.08854458077333
$ws.Cells.Item(3, 18).Value = 153.79690122696
$ws.Cells.Item(3, 19).Value = 0.1455965033634158
$ws.Cells.Item(3, 20).Value = 0.1455965033634159

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt9a"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.217872
$ws.Cells.Item(4, 8).Value = 3.653616
$ws.Cells.Item(4, 9).Value = 0.4812552202382371
$ws.Cells.Item(4, 10).Value = 0.4812552202382372
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.2022123333333333
$ws.Cells.Item(4, 14).Value = 0.606637
$ws.Cells.Item(4, 15).Value = 0.004359931864156574
$ws.Cells.Item(4, 16).Value = 0.004359931864156574
$ws.Cells.Item(4, 17).Value = 0.2462687388213333
$ws.Cells.Item(4, 18).Value = 2.216418649392
$ws.Cells.Item(4, 19).Value = 0.00209823996950838
$ws.Cells.Item(4, 20).Value = 0.00209823996950838

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt9a"
$ws.Cells.Item(5, 3).Value = "Fzd4"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.217872
$ws.Cells.Item(5, 8).Value = 3.653616
$ws.Cells.Item(5, 9).Value = 0.4812552202382371
$ws.Cells.Item(5, 10).Value = 0.4812552202382372
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.541567000000001
$ws.Cells.Item(5, 14).Value = 16.624701
$ws.Cells.Item(5, 15).Value = 0.1194825960532834
$ws.Cells.Item(5, 16).Value = 0.1194825960532834
$ws.Cells.Item(5, 17).Value = 6.748919285424001
$ws.Cells.Item(5, 18).Value = 60.74027356881601
$ws.Cells.Item(5, 19).Value = 0.05750162307825922
$ws.Cells.Item(5, 20).Value = 0.05750162307825923

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt9a"
$ws.Cells.Item(6, 3).Value = "Fzd4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.02442
$ws.Cells.Item(6, 8).Value = 3.07326
$ws.Cells.Item(6, 9).Value = 0.4048105816674124
$ws.Cells.Item(6, 10).Value = 0.4048105816674124
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 26.60444266666667
$ws.Cells.Item(6, 14).Value = 79.813328
$ws.Cells.Item(6, 15).Value = 0.5736225649467147
$ws.Cells.Item(6, 16).Value = 0.5736225649467147
$ws.Cells.Item(6, 17).Value = 27.25412315658667
$ws.Cells.Item(6, 18).Value = 245.28710840928
$ws.Cells.Item(6, 19).Value = 0.2322084841736326
$ws.Cells.Item(6, 20).Value = 0.2322084841736326

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt9a"
$ws.Cells.Item(7, 3).Value = "Fzd4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.02442
$ws.Cells.Item(7, 8).Value = 3.07326
$ws.Cells.Item(7, 9).Value = 0.4048105816674124
$ws.Cells.Item(7, 10).Value = 0.4048105816674124
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 14.03147833333333
$ws.Cells.Item(7, 14).Value = 42.094435
$ws.Cells.Item(7, 15).Value = 0.3025349071358453
$ws.Cells.Item(7, 16).Value = 0.3025349071358453
$ws.Cells.Item(7, 17).Value = 14.37412703423333
$ws.Cells.Item(7, 18).Value = 129.3671433081
$ws.Cells.Item(7, 19).Value = 0.1224693317323581
$ws.Cells.Item(7, 20).Value = 0.1224693317323581

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt9a"
$ws.Cells.Item(8, 3).Value = "Fzd4"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.02442
$ws.Cells.Item(8, 8).Value = 3.07326
$ws.Cells.Item(8, 9).Value = 0.4048105816674124
$ws.Cells.Item(8, 10).Value = 0.4048105816674124
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.2022123333333333
$ws.Cells.Item(8, 14).Value = 0.606637
$ws.Cells.Item(8, 15).Value = 0.004359931864156574
$ws.Cells.Item(8, 16).Value = 0.004359931864156574
$ws.Cells.Item(8, 17).Value = 0.2071503585133333
$ws.Cells.Item(8, 18).Value = 1.86435322662
$ws.Cells.Item(8, 19).Value = 0.001764946553959508
$ws.Cells.Item(8, 20).Value = 0.001764946553959509

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt9a"
$ws.Cells.Item(9, 3).Value = "Fzd4"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.02442
$ws.Cells.Item(9, 8).Value = 3.07326
$ws.Cells.Item(9, 9).Value = 0.4048105816674124
$ws.Cells.Item(9, 10).Value = 0.4048105816674124
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 5.541567000000001
$ws.Cells.Item(9, 14).Value = 16.624701
$ws.Cells.Item(9, 15).Value = 0.1194825960532834
$ws.Cells.Item(9, 16).Value = 0.1194825960532834
$ws.Cells.Item(9, 17).Value = 5.67689206614
$ws.Cells.Item(9, 18).Value = 51.09202859526
$ws.Cells.Item(9, 19).Value = 0.04836781920746212
$ws.Cells.Item(9, 20).Value = 0.04836781920746213

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Wnt9a"
$ws.Cells.Item(10, 3).Value = "Fzd4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.2883236666666666
$ws.Cells.Item(10, 8).Value = 0.8649709999999999
$ws.Cells.Item(10, 9).Value = 0.1139341980943504
$ws.Cells.Item(10, 10).Value = 0.1139341980943504
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 26.60444266666667
$ws.Cells.Item(10, 14).Value = 79.813328
$ws.Cells.Item(10, 15).Value = 0.5736225649467147
$ws.Cells.Item(10, 16).Value = 0.5736225649467147
$ws.Cells.Item(10, 17).Value = 7.670690459276444
$ws.Cells.Item(10, 18).Value = 69.03621413348799
$ws.Cells.Item(10, 19).Value = 0.06535522694602838
$ws.Cells.Item(10, 20).Value = 0.06535522694602838

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Wnt9a"
$ws.Cells.Item(11, 3).Value = "Fzd4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.2883236666666666
$ws.Cells.Item(11, 8).Value = 0.8649709999999999
$ws.Cells.Item(11, 9).Value = 0.1139341980943504
$ws.Cells.Item(11, 10).Value = 0.1139341980943504
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 14.03147833333333
$ws.Cells.Item(11, 14).Value = 42.094435
$ws.Cells.Item(11, 15).Value = 0.3025349071358453
$ws.Cells.Item(11, 16).Value = 0.3025349071358453
$ws.Cells.Item(11, 17).Value = 4.045607281820556
$ws.Cells.Item(11, 18).Value = 36.410465536385
$ws.Cells.Item(11, 19).Value = 0.0344690720400713
$ws.Cells.Item(11, 20).Value = 0.03446907204007131

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Wnt9a"
$ws.Cells.Item(12, 3).Value = "Fzd4"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.2883236666666666
$ws.Cells.Item(12, 8).Value = 0.8649709999999999
$ws.Cells.Item(12, 9).Value = 0.1139341980943504
$ws.Cells.Item(12, 10).Value = 0.1139341980943504
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.2022123333333333
$ws.Cells.Item(12, 14).Value = 0.606637
$ws.Cells.Item(12, 15).Value = 0.004359931864156574
$ws.Cells.Item(12, 16).Value = 0.004359931864156574
$ws.Cells.Item(12, 17).Value = 0.05830260139188888
$ws.Cells.Item(12, 18).Value = 0.5247234125269999
$ws.Cells.Item(12, 19).Value = 0.0004967453406886856
$ws.Cells.Item(12, 20).Value = 0.0004967453406886857

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Wnt9a"
$ws.Cells.Item(13, 3).Value = "Fzd4"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.2883236666666666
$ws.Cells.Item(13, 8).Value = 0.8649709999999999
$ws.Cells.Item(13, 9).Value = 0.1139341980943504
$ws.Cells.Item(13, 10).Value = 0.1139341980943504
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 5.541567000000001
$ws.Cells.Item(13, 14).Value = 16.624701
$ws.Cells.Item(13, 15).Value = 0.1194825960532834
$ws.Cells.Item(13, 16).Value = 0.1194825960532834
$ws.Cells.Item(13, 17).Value = 1.597764916519
$ws.Cells.Item(13, 18).Value = 14.379884248671
$ws.Cells.Item(13, 19).Value = 0.01361315376756204
$ws.Cells.Item(13, 20).Value = 0.01361315376756204

# Remove the now-obsolete rows (previously rows 14-16, for the "Resolving-Mac" target cluster)
$ws.Rows(14).Resize(3).Delete() | Out-Null
